$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 473.875
$ws.Range("I5").Value = 470.14285
$ws.Range("K5").Value = 470.14285
$ws.Range("M5").Value = -355.14285
$ws.Range("H80").Value = 856.2727
$ws.Range("I80").Value = 162
$ws.Range("J80").Value = 1116.625
$ws.Range("K80").Value = 486
$ws.Range("L80").Value = 3349.875
$ws.Range("M80").Value = 512
$ws.Range("N80").Value = -5345.875
$ws.Range("H83").Value = 856.2727
$ws.Range("I83").Value = 162
$ws.Range("J83").Value = 1116.625
$ws.Range("K83").Value = 1458
$ws.Range("L83").Value = 10049.625
$ws.Range("M83").Value = 3534
$ws.Range("N83").Value = -20033.625
$ws.Range("H98").Value = 1226.2174
$ws.Range("I98").Value = 1266.5814
$ws.Range("K98").Value = 1266.5814
$ws.Range("M98").Value = 231.4186
$ws.Range("H122").Value = 1226.2174
$ws.Range("I122").Value = 1266.5814
$ws.Range("K122").Value = 3799.7442
$ws.Range("M122").Value = -1349.7442
$ws.Range("H138").Value = 5884773
$ws.Range("J138").Value = 8775027
$ws.Range("L138").Value = 26325081
$ws.Range("N138").Value = -26335361

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7460.8945
$ws.Range("I32").Value = 4887.7646
$ws.Range("K32").Value = 4887.7646
$ws.Range("M32").Value = -4600.7646
$ws.Range("H45").Value = 14937.125
$ws.Range("I45").Value = 16713.857
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 16713.857
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -16336.857
$ws.Range("N45").Value = -3254
$ws.Range("H122").Value = 1791.875
$ws.Range("I122").Value = 1500.4445
$ws.Range("J122").Value = 2666.1667
$ws.Range("K122").Value = 4501.333500000001
$ws.Range("L122").Value = 7998.500100000001
$ws.Range("M122").Value = -2051.333500000001
$ws.Range("N122").Value = -12898.5001
$ws.Range("H139").Value = 97499.39999999999
$ws.Range("J139").Value = 97499.39999999999
$ws.Range("L139").Value = 97499.39999999999
$ws.Range("N139").Value = -107779.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 286.90475
$ws.Range("I7").Value = 336.3
$ws.Range("J7").Value = 242
$ws.Range("K7").Value = 336.3
$ws.Range("L7").Value = 242
$ws.Range("M7").Value = -223.3
$ws.Range("N7").Value = -468
$ws.Range("H22").Value = 187.33333
$ws.Range("I22").Value = 187.33333
$ws.Range("K22").Value = 187.33333
$ws.Range("M22").Value = 162.66667
$ws.Range("H31").Value = 52922.6
$ws.Range("I31").Value = 86015.414
$ws.Range("K31").Value = 86015.414
$ws.Range("M31").Value = -85720.414
$ws.Range("H34").Value = 52922.6
$ws.Range("I34").Value = 86015.414
$ws.Range("K34").Value = 86015.414
$ws.Range("M34").Value = -85813.414
$ws.Range("H99").Value = 4098.6665
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4098.6665
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -7094.6665
$ws.Range("H126").Value = 4098.6665
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4098.6665
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -17235.9995
$ws.Range("H134").Value = 9430.102000000001
$ws.Range("I134").Value = 4228.976
$ws.Range("J134").Value = 40636.855
$ws.Range("K134").Value = 12686.928
$ws.Range("L134").Value = 121910.565
$ws.Range("M134").Value = -10151.928
$ws.Range("N134").Value = -126980.565

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 838.125
$ws.Range("I113").Value = 803.8889
$ws.Range("J113").Value = 858.6667
$ws.Range("K113").Value = 2411.6667
$ws.Range("L113").Value = 2576.0001
$ws.Range("M113").Value = -241.6667000000002
$ws.Range("N113").Value = -6916.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5700000
$ws.Range("J11").Value = 50000
$ws.Range("L11").Value = 50000
$ws.Range("N11").Value = -50278
$ws.Range("H18").Value = 55565550
$ws.Range("I18").Value = 55565550
$ws.Range("K18").Value = 55565550
$ws.Range("M18").Value = -55565257
$ws.Range("H22").Value = 33000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 33000
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -34058
$ws.Range("H24").Value = 44999.145
$ws.Range("J24").Value = 44999.145
$ws.Range("L24").Value = 44999.145
$ws.Range("N24").Value = -45345.145
$ws.Range("H122").Value = 2984.862
$ws.Range("I122").Value = 2524.1738
$ws.Range("K122").Value = 7572.5214
$ws.Range("M122").Value = -5122.5214
$ws.Range("H132").Value = 3160.0908
$ws.Range("I132").Value = 2911.6843
$ws.Range("J132").Value = 4733.3335
$ws.Range("K132").Value = 8735.052899999999
$ws.Range("L132").Value = 14200.0005
$ws.Range("M132").Value = -6205.052899999999
$ws.Range("N132").Value = -19260.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1688.6666
$ws.Range("I16").Value = 1752.1923
$ws.Range("K16").Value = 1752.1923
$ws.Range("M16").Value = -1582.1923
$ws.Range("H22").Value = 3258.4546
$ws.Range("I22").Value = 3141
$ws.Range("J22").Value = 3399.4
$ws.Range("K22").Value = 3141
$ws.Range("L22").Value = 3399.4
$ws.Range("M22").Value = -2846
$ws.Range("N22").Value = -3989.4
$ws.Range("H23").Value = 50000000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 50000000
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -50000460
$ws.Range("H25").Value = 2212.121
$ws.Range("I25").Value = 2212.121
$ws.Range("K25").Value = 2212.121
$ws.Range("M25").Value = -1982.121
$ws.Range("H27").Value = 3258.4546
$ws.Range("I27").Value = 3141
$ws.Range("J27").Value = 3399.4
$ws.Range("K27").Value = 3141
$ws.Range("L27").Value = 3399.4
$ws.Range("M27").Value = -3034
$ws.Range("N27").Value = -3613.4
$ws.Range("H40").Value = 4636.0303
$ws.Range("I40").Value = 3843.5715
$ws.Range("K40").Value = 3843.5715
$ws.Range("M40").Value = -3707.5715
$ws.Range("H68").Value = 3069
$ws.Range("I68").Value = 2962.4443
$ws.Range("K68").Value = 2962.4443
$ws.Range("M68").Value = -2213.4443
$ws.Range("H71").Value = 3069
$ws.Range("I71").Value = 2962.4443
$ws.Range("K71").Value = 14812.2215
$ws.Range("M71").Value = -11068.2215
$ws.Range("H122").Value = 4560.811
$ws.Range("I122").Value = 3708.2903
$ws.Range("K122").Value = 11124.8709
$ws.Range("M122").Value = -8674.8709
$ws.Range("H132").Value = 6807.1816
$ws.Range("I132").Value = 4983.1665
$ws.Range("K132").Value = 14949.4995
$ws.Range("M132").Value = -12419.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 20000000
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 20000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 20000
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -21802
$ws.Range("H136").Value = 2131.054
$ws.Range("I136").Value = 1966
$ws.Range("K136").Value = 5898
$ws.Range("M136").Value = -3348
